$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 3.600374698638916
$ws.Range("B1").Value = 2.693195104598999
$ws.Range("C1").Value = 2.26577615737915
$ws.Range("D1").Value = 2.355586528778076
$ws.Range("E1").Value = 2.665621519088745
